# "commit final, trabajo terminado"
#
# Row 12 ("Edición de Ropa y Zapatos" test case) is finished/approved:
#  - The "Resultados actuales" (E12) text is updated to match the
#    "Resultados esperados" (D12) text (the test now behaves as expected,
#    instead of the old "duplicates the product" bug description).
#  - The "Aprobado?" column (F12) is filled in with "Si" now that the
#    test case passed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$expected = $ws.Range("D12").Value()
$ws.Range("E12").Value = $expected
$ws.Range("F12").Value = "Si"

# Leave the selection where the work finished, on the result cell of row 11.
[void]$ws.Range("E11").Select()
